$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The updated LR-pair results only contain 6 data rows (rows 2-7) instead of 9,
# because the "ECs" sending-cluster rows were dropped. Remove the trailing rows first.
$ws.Rows("8:10").Delete()

# New data for rows 2-7: Sending cluster, Ligand, Receptor, Target cluster,
# followed by the 16 recomputed TPM-derived metric columns (E:T).
$data = @(
    @("FAPs", "Slit1", "Sdc1", "ECs", 2, 0.6666666666666666, 0.2055996666666667, 0.616799, 0.9059768423248155, 0.9059768423248156, 3, 1, 0.8213140000000001, 2.463942, 0.06824749762056036, 0.06824749762056037, 0.1688618846286667, 1.519756961658, 0.06183065239084563, 0.06183065239084565),
    @("FAPs", "Slit1", "Sdc1", "FAPs", 2, 0.6666666666666666, 0.2055996666666667, 0.616799, 0.9059768423248155, 0.9059768423248156, 3, 1, 2.427350333333333, 7.282051, 0.2017018900182306, 0.2017018900182306, 0.4990624194165555, 4.491561774749, 0.1827372414096638, 0.1827372414096638),
    @("FAPs", "Slit1", "Sdc1", "MuSCs", 2, 0.6666666666666666, 0.2055996666666667, 0.616799, 0.9059768423248155, 0.9059768423248156, 3, 1, 8.785681666666667, 26.357045, 0.730050612361209, 0.730050612361209, 1.806333222106111, 16.256998998955, 0.6614089485243061, 0.6614089485243061),
    @("MuSCs", "Slit1", "Sdc1", "ECs", 2, 0.6666666666666666, 0.02133733333333333, 0.064012, 0.09402315767518445, 0.09402315767518446, 3, 1, 0.8213140000000001, 2.463942, 0.06824749762056036, 0.06824749762056037, 0.01752465058933333, 0.157721855304, 0.006416845229714722, 0.006416845229714724),
    @("MuSCs", "Slit1", "Sdc1", "FAPs", 2, 0.6666666666666666, 0.02133733333333333, 0.064012, 0.09402315767518445, 0.09402315767518446, 3, 1, 2.427350333333333, 7.282051, 0.2017018900182306, 0.2017018900182306, 0.05179318317911111, 0.466138648612, 0.01896464860856681, 0.01896464860856681),
    @("MuSCs", "Slit1", "Sdc1", "MuSCs", 2, 0.6666666666666666, 0.02133733333333333, 0.064012, 0.09402315767518445, 0.09402315767518446, 3, 1, 8.785681666666667, 26.357045, 0.730050612361209, 0.730050612361209, 0.1874630182822222, 1.68716716454, 0.06864166383690291, 0.06864166383690293),
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$r").Value = $row[$j]
    }
}
